# Weekly roll-forward update for the Mango / Terminal La Palmera de La Serena sheet.
# A new week's price group (3 rows: Especial/Primera/Segunda) is inserted at the
# top of the data block (rows 741-743) and every existing data row below it
# shifts down by 3 rows, growing the table by 3 rows (792 -> 795).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the existing data block (rows 741-792) down by 3 rows (-> 744-795).
#    Value2 bulk copy preserves every column's content (dates stay as serials).
$srcValues = $ws.Range("A741:T792").Value2
$ws.Range("A744:T795").Value2 = $srcValues

# 2) Extend the date number-format (column D) down onto the newly created rows
#    (793-795) by copying the existing date-formatted cells' formatting.
$ws.Range("D741:D792").Copy()
$ws.Range("D744:D795").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Write the new top group (rows 741-743) with this week's figures.
$row741 = New-Object 'object[,]' 1,20
$row741[0,0]  = 8
$row741[0,1]  = "Terminal La Palmera de La Serena"
$row741[0,2]  = "Coquimbo"
$row741[0,3]  = 44746
$row741[0,4]  = 4
$row741[0,5]  = "Fruta"
$row741[0,6]  = 100108
$row741[0,7]  = "Tropicales y subtropicales"
$row741[0,8]  = 100108002
$row741[0,9]  = "Mango"
$row741[0,10] = "Sin especificar"
$row741[0,11] = "Especial"
$row741[0,12] = 512
$row741[0,13] = 7500
$row741[0,14] = 8000
$row741[0,15] = 7750
$row741[0,16] = "`$/bandeja 4 kilos"
$row741[0,17] = "Brasil"
$row741[0,18] = 1938
$row741[0,19] = 4
$ws.Range("A741:T741").Value2 = $row741

$row742 = New-Object 'object[,]' 1,20
$row742[0,0]  = 8
$row742[0,1]  = "Terminal La Palmera de La Serena"
$row742[0,2]  = "Coquimbo"
$row742[0,3]  = 44746
$row742[0,4]  = 4
$row742[0,5]  = "Fruta"
$row742[0,6]  = 100108
$row742[0,7]  = "Tropicales y subtropicales"
$row742[0,8]  = 100108002
$row742[0,9]  = "Mango"
$row742[0,10] = "Sin especificar"
$row742[0,11] = "Primera"
$row742[0,12] = 515
$row742[0,13] = 7500
$row742[0,14] = 8000
$row742[0,15] = 7749
$row742[0,16] = "`$/bandeja 4 kilos"
$row742[0,17] = "Brasil"
$row742[0,18] = 1937
$row742[0,19] = 4
$ws.Range("A742:T742").Value2 = $row742

$row743 = New-Object 'object[,]' 1,20
$row743[0,0]  = 8
$row743[0,1]  = "Terminal La Palmera de La Serena"
$row743[0,2]  = "Coquimbo"
$row743[0,3]  = 44746
$row743[0,4]  = 4
$row743[0,5]  = "Fruta"
$row743[0,6]  = 100108
$row743[0,7]  = "Tropicales y subtropicales"
$row743[0,8]  = 100108002
$row743[0,9]  = "Mango"
$row743[0,10] = "Sin especificar"
$row743[0,11] = "Segunda"
$row743[0,12] = 512
$row743[0,13] = 7500
$row743[0,14] = 8000
$row743[0,15] = 7750
$row743[0,16] = "`$/bandeja 4 kilos"
$row743[0,17] = "Brasil"
$row743[0,18] = 1938
$row743[0,19] = 4
$ws.Range("A743:T743").Value2 = $row743
